# ZBP_13_dusevni_zdravi.xlsx — "Add files via upload" update
#
# Adds a new survey wave, "12. 10. 2021", as the last column on both
# worksheets ("data" = percentages, "pocetR" = sample sizes), and bumps
# the "aktualizace" (last-updated) date in both footer/title cells from
# "6. 10. 2021" to "20. 10. 2021".

$wb  = $excel.ActiveWorkbook
$wsData   = $wb.Worksheets.Item("data")
$wsPocetR = $wb.Worksheets.Item("pocetR")

# ---------------------------------------------------------------------
# Sheet "data": new column AI, header + 22 rows of percentages
# ---------------------------------------------------------------------

# Clone the formatting of the last existing header cell (AH1 -> AI1) so
# the new header keeps the same bold/centered/bordered look, then set
# its text.
$wsData.Range("AH1").Copy()
$wsData.Range("AI1").PasteSpecial(-4122)  # xlPasteFormats
$wsData.Range("AI1").Value = "12. 10. 2021"

$dataValues = [ordered]@{
    2  = 0.1
    3  = 0.09
    4  = 0.11
    5  = 0.06
    6  = 0.1
    7  = 0.15
    8  = 0.1
    9  = 0.2
    10 = 0.1
    11 = 0.08
    12 = 0.11
    13 = 0.08
    14 = 0.17
    15 = 0.1
    16 = 0.09
    17 = 0.1
    18 = 0.12
    19 = 0.07
    20 = 0.09
    21 = 0.06
    22 = 0.07
    23 = 0.18
}

foreach ($row in $dataValues.Keys) {
    $wsData.Range("AI$row").Value = $dataValues[$row]
}

# Footer title cell: bump the "aktualizace" date
$wsData.Range("A24").Value = "Život během pandemie, Duševní zdraví, % respondentů celkově a ve skupinách, aktualizace 20. 10. 2021"

# ---------------------------------------------------------------------
# Sheet "pocetR": new column AH, header + 22 rows of sample sizes
# ---------------------------------------------------------------------

$wsPocetR.Range("AG1").Copy()
$wsPocetR.Range("AH1").PasteSpecial(-4122)  # xlPasteFormats
$wsPocetR.Range("AH1").Value = "12. 10. 2021"

$pocetRValues = [ordered]@{
    2  = 1836
    3  = 895
    4  = 941
    5  = 239
    6  = 655
    7  = 281
    8  = 661
    9  = 158
    10 = 296
    11 = 348
    12 = 322
    13 = 712
    14 = 187
    15 = 352
    16 = 1297
    17 = 193
    18 = 666
    19 = 611
    20 = 245
    21 = 522
    22 = 819
    23 = 495
}

foreach ($row in $pocetRValues.Keys) {
    $wsPocetR.Range("AH$row").Value = $pocetRValues[$row]
}

# Footer title cell: bump the "aktualizace" date
$wsPocetR.Range("A24").Value = "Život během pandemie, Duševní zdraví, velikost dotázaného souboru celkově a ve skupinách, aktualizace 20. 10. 2021"

Write-Output "done"
